# Change the table style (table gallery pick) on the "Sources of finance"
# table from the deck's custom "Table_0" style to the built-in
# "Medium Style 2 - Accent 1" gallery style ({08846E95-98BA-494B-A4F6-413BB014342C}).
#
# The table lives on slide 6 (the only slide/shape in the deck with a table).
# We locate it defensively (scan every slide/shape) instead of hard-coding
# indices, so the script keeps working even if slide/shape ordering differs.

$p = $ppt.ActivePresentation

$targetStyleId = "{08846E95-98BA-494B-A4F6-413BB014342C}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($targetStyleId)
        }
    }
}
